# Apply updated cryptocurrency market data (price & 1h volume change)
# to the active worksheet, preserving each target cell as plain text
# (matching the inline-string cell type used throughout the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    # Force text interpretation so numeric-looking strings (e.g. '676.32')
    # are not silently coerced into Excel numbers.
    $cell.NumberFormat = '@'
    $cell.Value = $Text
    # Drop the quote-prefix/text style that was just applied so the cell's
    # formatting matches the rest of the (unstyled) data cells.
    $cell.Style = 'Normal'
}

Set-TextCell 'D2' '69.290.95'
Set-TextCell 'E2' '  -0.07%  '
Set-TextCell 'D3' '3.674.05'
Set-TextCell 'E3' '  -0.31%  '
Set-TextCell 'E4' '  +0.02%  '
Set-TextCell 'D5' '676.32'
Set-TextCell 'E5' '  -0.77%  '
Set-TextCell 'D6' '158.51'
Set-TextCell 'E6' '  -2.34%  '
Set-TextCell 'E7' '  -0.05%  '
Set-TextCell 'E8' '  -1.33%  '
Set-TextCell 'E9' '  -1.34%  '
Set-TextCell 'D10' '6.93'
Set-TextCell 'E10' '  -5.40%  '
Set-TextCell 'E11' '  -2.64%  '
Set-TextCell 'E12' '  -3.04%  '
Set-TextCell 'D13' '4.294.45'
Set-TextCell 'E13' '  -0.34%  '
Set-TextCell 'D14' '32.37'
Set-TextCell 'E14' '  -3.76%  '
Set-TextCell 'D15' '3.679.85'
Set-TextCell 'E15' '  -0.19%  '
Set-TextCell 'D16' '69.238.60'
Set-TextCell 'E16' '  -0.22%  '
Set-TextCell 'E17' '  +1.56%  '
Set-TextCell 'D18' '16.05'
Set-TextCell 'E18' '  -1.77%  '
Set-TextCell 'D19' '6.44'
Set-TextCell 'E19' '  -2.74%  '
Set-TextCell 'D20' '467.36'
Set-TextCell 'E20' '  -2.95%  '
Set-TextCell 'D21' '10.01'
Set-TextCell 'E21' '  +0.99%  '
Set-TextCell 'D22' '0.649'
Set-TextCell 'E22' '  -2.80%  '
Set-TextCell 'D23' '79.74'
Set-TextCell 'E23' '  -0.59%  '
Set-TextCell 'D24' '3.819.19'
Set-TextCell 'E24' '  -0.34%  '
Set-TextCell 'E25' '  +0.02%  '
Set-TextCell 'E26' '  -6.28%  '
Set-TextCell 'E27' '  -5.17%  '
Set-TextCell 'D28' '9.07'
Set-TextCell 'E28' '  -4.77%  '
Set-TextCell 'D29' '2.68'
Set-TextCell 'E29' '  -1.02%  '
Set-TextCell 'D30' '1.75'
Set-TextCell 'E30' '  -4.82%  '
Set-TextCell 'D31' '6.61'
Set-TextCell 'E31' '  -3.57%  '
Set-TextCell 'B32' 'ImmutableX'
Set-TextCell 'C32' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D32' '1.99'
Set-TextCell 'E32' '  -4.69%  '
Set-TextCell 'B33' 'Binance-PegBSC-USD'
Set-TextCell 'C33' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell 'D33' '0.999'
Set-TextCell 'E33' '  +0.03%  '
Set-TextCell 'D34' '26.89'
Set-TextCell 'E34' '  -0.69%  '
Set-TextCell 'D35' '3.667.50'
Set-TextCell 'E35' '  +0.29%  '
Set-TextCell 'D36' '0.162'
Set-TextCell 'E36' '  -4.57%  '
Set-TextCell 'D37' '8.19'
Set-TextCell 'E37' '  -3.35%  '
Set-TextCell 'D38' '6.20'
Set-TextCell 'E38' '  -2.00%  '
Set-TextCell 'E39' '  +0.01%  '
Set-TextCell 'E40' '  -0.13%  '
Set-TextCell 'D41' '2.23'
Set-TextCell 'E41' '  -1.34%  '
Set-TextCell 'B42' 'Hedera'
Set-TextCell 'C42' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D42' '0.0901'
Set-TextCell 'E42' '  -3.83%  '
Set-TextCell 'B43' 'Monero'
Set-TextCell 'C43' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D43' '174.34'
Set-TextCell 'E43' '  +7.75%  '
Set-TextCell 'D44' '0.940'
Set-TextCell 'E44' '  -1.51%  '
Set-TextCell 'D45' '47.56'
Set-TextCell 'E45' '  -1.70%  '
Set-TextCell 'D46' '28.19'
Set-TextCell 'E46' '  -6.37%  '
Set-TextCell 'D47' '2.70'
Set-TextCell 'E47' '  -4.73%  '
Set-TextCell 'D48' '0.000278'
Set-TextCell 'E48' '  -3.20%  '
Set-TextCell 'E49' '  -4.51%  '
Set-TextCell 'E50' '  -3.69%  '
Set-TextCell 'E51' '  -3.09%  '
